# ----------------------------------------------------------------------------
# Applies the "B0BHM59TQB_po_data.xlsx" edit:
#   1. Rename "Weekly Quantity" header B1  -> "Weekly_PO_Qty"
#   2. Rename "Monthly Trend" header B1    -> "Monthly_PO_Qty"
#   3. Add a new "PO Forecast" worksheet at the end with forecast data
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# 1) Header rename on "Weekly Quantity"
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Header rename on "Monthly Trend"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Build the new "PO Forecast" sheet by duplicating "Weekly Quantity" (so it
#    inherits identical sheet-level properties: outline settings, page setup,
#    margins, etc.), placing the duplicate after the last existing tab, then
#    clearing its old contents (keeping the style of the cells in place) and
#    writing the forecast data into it.
$wsWeekly.Copy($null, $wb.Worksheets($wb.Worksheets.Count))
$wsForecast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast.Name = "PO Forecast"
$wsForecast.Cells.ClearContents()

# Extend the inherited header style (currently only on A1:B1) across to D1,
# and the inherited date-column style (currently only on A2:A7) down to A15.
$wsForecast.Range("A1:B1").Copy($wsForecast.Range("A1:D1"))
$wsForecast.Range("A2:A7").Copy($wsForecast.Range("A8:A15"))

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast data rows
$wsForecast.Range("A2").Value = 45011.99999999999
$wsForecast.Range("B2").Value = 50
$wsForecast.Range("C2").Value = 49.99999994175453
$wsForecast.Range("D2").Value = 50.0000000628256
$wsForecast.Range("A3").Value = 45081.99999999999
$wsForecast.Range("B3").Value = 50
$wsForecast.Range("C3").Value = 49.99999993360198
$wsForecast.Range("D3").Value = 50.00000006274355
$wsForecast.Range("A4").Value = 45088.99999999999
$wsForecast.Range("B4").Value = 50
$wsForecast.Range("C4").Value = 49.99999993643041
$wsForecast.Range("D4").Value = 50.00000005999608
$wsForecast.Range("A5").Value = 45102.99999999999
$wsForecast.Range("B5").Value = 50
$wsForecast.Range("C5").Value = 49.99999993341022
$wsForecast.Range("D5").Value = 50.00000006160847
$wsForecast.Range("A6").Value = 45116.99999999999
$wsForecast.Range("B6").Value = 50
$wsForecast.Range("C6").Value = 49.99999993776879
$wsForecast.Range("D6").Value = 50.00000006533416
$wsForecast.Range("A7").Value = 45144.99999999999
$wsForecast.Range("B7").Value = 50
$wsForecast.Range("C7").Value = 49.99999993224107
$wsForecast.Range("D7").Value = 50.00000006504066
$wsForecast.Range("A8").Value = 45151.99999999999
$wsForecast.Range("B8").Value = 50
$wsForecast.Range("C8").Value = 49.99999993528606
$wsForecast.Range("D8").Value = 50.00000006382224
$wsForecast.Range("A9").Value = 45158.99999999999
$wsForecast.Range("B9").Value = 50
$wsForecast.Range("C9").Value = 49.9999999349939
$wsForecast.Range("D9").Value = 50.000000068103
$wsForecast.Range("A10").Value = 45165.99999999999
$wsForecast.Range("B10").Value = 50
$wsForecast.Range("C10").Value = 49.99999992495491
$wsForecast.Range("D10").Value = 50.00000007649221
$wsForecast.Range("A11").Value = 45172.99999999999
$wsForecast.Range("B11").Value = 50
$wsForecast.Range("C11").Value = 49.9999999145096
$wsForecast.Range("D11").Value = 50.00000008932096
$wsForecast.Range("A12").Value = 45179.99999999999
$wsForecast.Range("B12").Value = 50
$wsForecast.Range("C12").Value = 49.99999989835559
$wsForecast.Range("D12").Value = 50.00000010376758
$wsForecast.Range("A13").Value = 45186.99999999999
$wsForecast.Range("B13").Value = 50
$wsForecast.Range("C13").Value = 49.99999986735479
$wsForecast.Range("D13").Value = 50.00000011961817
$wsForecast.Range("A14").Value = 45193.99999999999
$wsForecast.Range("B14").Value = 50
$wsForecast.Range("C14").Value = 49.99999982075617
$wsForecast.Range("D14").Value = 50.0000001537988
$wsForecast.Range("A15").Value = 45200.99999999999
$wsForecast.Range("B15").Value = 50
$wsForecast.Range("C15").Value = 49.99999979201425
$wsForecast.Range("D15").Value = 50.00000019514172
